$d = $word.ActiveDocument

# 1. Fix typo: "vælger Linær" -> "vælger lineær" (main success scenario step about
#    choosing the depreciation method).
$d.Content.Find.Execute("vælger Linær", $true, $false, $false, $false, $false,
                         $true, 1, $false, "vælger lineær", 2)

# 2. Miscellaneous row: replace the placeholder "N/A" with a real comment about
#    the statutory max amount for "straksafskrivning" (2b.3a).
#    Scope the replacement strictly to that table cell (re-anchor the Range via
#    Document.Range(start,end) so Find doesn't wander off to the first "N/A" in
#    the whole document), and replace only the first (only) match in range.
$miscCell = $d.Tables(1).Cell(14, 2)
$cellRange = $d.Range($miscCell.Range.Start, $miscCell.Range.End)
$cellRange.Find.Execute("N/A", $true, $false, $false, $false, $false,
                         $true, 0, $false,
                         "Max-beløbet for straksafskrivning (2b.3a) reguleres lovgivningsmæssigt en gang om året.",
                         1)
